$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for the two new columns, matching style of existing header (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data for columns I (I0) and J (IF), rows 2 through 44
$iValues = @(6, 4, 6, 6, 8, 6, 10, 6, 6, 9, 8, 6, 8, 10, 5, 8, 6, 7, 8, 7, 9, 7, 8, 7, 6, 8, 6, 9, 8, 8, 7, 7, 7, 8, 6, 7, 8, 6, 6, 5, 5, 7, 4)
$jValues = @(6, 5, 6, 6, 9, 6, 10, 7, 7, 9, 8, 7, 8, 11, 5, 8, 6, 7, 8, 8, 9, 7, 9, 7, 7, 8, 6, 9, 8, 8, 8, 8, 8, 9, 8, 8, 8, 6, 6, 6, 6, 7, 5)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
